$wb = $excel.ActiveWorkbook
$wsCbc = $wb.Worksheets.Item("CBC")

# --- Update the CBC source-data table (rows 5-12, cols A-H) ---
# Columns B, C, D switch from text ("100.0","100.0","60.0") to real numbers (100,100,60)
# Column E (SKU), F (Loaded2) and G (Column #) are refactored with two extra rows.

$data = @(
    @(7,  "20842-ESCALOPIN LOMO ADOB. 300 GR. B", 1),
    @(8,  "21277-LOMO SIEMPRE TIERNO 300 GR. BA", 2),
    @(7,  "21277-LOMO SIEMPRE TIERNO 300 GR. BA", 3),
    @(7,  "21277-LOMO SIEMPRE TIERNO 300 GR. BA", 4),
    @(7,  "20842-ESCALOPIN LOMO ADOB. 300 GR. B", 5),
    @(8,  "21277-LOMO SIEMPRE TIERNO 300 GR. BA", 6),
    @(8,  "21277-LOMO SIEMPRE TIERNO 300 GR. BA", 7),
    @(8,  "21277-LOMO SIEMPRE TIERNO 300 GR. BA", 8)
)

$row = 5
foreach ($item in $data) {
    $wsCbc.Cells.Item($row, 1).Value = "PALET-001(352-1984754-PO)"
    $wsCbc.Cells.Item($row, 2).Value = 100
    $wsCbc.Cells.Item($row, 3).Value = 100
    $wsCbc.Cells.Item($row, 4).Value = 60
    $wsCbc.Cells.Item($row, 5).Value = $item[1]
    $wsCbc.Cells.Item($row, 6).Value = $item[0]
    $wsCbc.Cells.Item($row, 7).Value = $item[2]
    $wsCbc.Cells.Item($row, 8).Value = 999999
    $row = $row + 1
}

# --- Refresh both pivot tables so their caches/records/output reflect the new data ---
$wsMain = $wb.Worksheets.Item("Comparativa CBC CubeMaster")
$wsMain.PivotTables(1).RefreshTable()
$wsMain.PivotTables(2).RefreshTable()

Write-Output "done"
